# Updates India Super League odds data: corrects results/odds for existing
# fixtures (rows 103-104) and appends newly played/scheduled fixtures
# (rows 105-107), per league database refresh on 07-03-2024.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 103: update existing match result / closing odds ---
$ws.Cells.Item(103, 8).Value = 2
$ws.Cells.Item(103, 9).Value = 1
$ws.Cells.Item(103, 10).Value = "H"
$ws.Cells.Item(103, 14).Value = 2.625
$ws.Cells.Item(103, 15).Value = 3.4
$ws.Cells.Item(103, 16).Value = 2.25
$ws.Cells.Item(103, 17).Value = 0
$ws.Cells.Item(103, 18).Value = 2.05
$ws.Cells.Item(103, 19).Value = 1.8
$ws.Cells.Item(103, 20).Value = 2.5
$ws.Cells.Item(103, 21).Value = 1.925
$ws.Cells.Item(103, 22).Value = 1.925
$ws.Cells.Item(103, 23).Value = 1.625
$ws.Cells.Item(103, 24).Value = -1
$ws.Cells.Item(103, 25).Value = -1
$ws.Cells.Item(103, 26).Value = 1.05
$ws.Cells.Item(103, 27).Value = -1
$ws.Cells.Item(103, 28).Value = 0.925
$ws.Cells.Item(103, 29).Value = -1

# --- Row 104: update existing match result / closing odds ---
$ws.Cells.Item(104, 8).Value = 2
$ws.Cells.Item(104, 9).Value = 2
$ws.Cells.Item(104, 10).Value = "D"
$ws.Cells.Item(104, 15).Value = 4.5
$ws.Cells.Item(104, 19).Value = 1.9
$ws.Cells.Item(104, 21).Value = 1.85
$ws.Cells.Item(104, 22).Value = 2
$ws.Cells.Item(104, 23).Value = -1
$ws.Cells.Item(104, 24).Value = 3.5
$ws.Cells.Item(104, 25).Value = -1
$ws.Cells.Item(104, 26).Value = 0.95
$ws.Cells.Item(104, 27).Value = -1
$ws.Cells.Item(104, 28).Value = 0.8500000000000001
$ws.Cells.Item(104, 29).Value = -1

# --- Row 105: new fixture ---
$ws.Range("A2").Copy()
$ws.Cells.Item(105, 1).PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Cells.Item(105, 5).PasteSpecial(-4122)
$ws.Cells.Item(105, 1).Value = 103
$ws.Cells.Item(105, 2).Value = 7751755
$ws.Cells.Item(105, 3).Value = "India Super League"
$ws.Cells.Item(105, 4).Value = "India Super League"
$ws.Cells.Item(105, 5).Value = 45357.45833333334
$ws.Cells.Item(105, 6).Value = "FC Goa"
$ws.Cells.Item(105, 7).Value = "East Bengal Club"
$ws.Cells.Item(105, 8).Value = 1
$ws.Cells.Item(105, 9).Value = 0
$ws.Cells.Item(105, 10).Value = "H"
$ws.Cells.Item(105, 11).Value = 1.6
$ws.Cells.Item(105, 12).Value = 4
$ws.Cells.Item(105, 13).Value = 5.25
$ws.Cells.Item(105, 14).Value = 1.6
$ws.Cells.Item(105, 15).Value = 4
$ws.Cells.Item(105, 16).Value = 4.75
$ws.Cells.Item(105, 17).Value = -1
$ws.Cells.Item(105, 18).Value = 2.05
$ws.Cells.Item(105, 19).Value = 1.8
$ws.Cells.Item(105, 20).Value = 2.75
$ws.Cells.Item(105, 21).Value = 1.875
$ws.Cells.Item(105, 22).Value = 1.975
$ws.Cells.Item(105, 23).Value = 0.6000000000000001
$ws.Cells.Item(105, 24).Value = -1
$ws.Cells.Item(105, 25).Value = -1
$ws.Cells.Item(105, 26).Value = 0
$ws.Cells.Item(105, 27).Value = -0.0
$ws.Cells.Item(105, 28).Value = -1
$ws.Cells.Item(105, 29).Value = 0.9750000000000001

# --- Row 106: new fixture ---
$ws.Range("A2").Copy()
$ws.Cells.Item(106, 1).PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Cells.Item(106, 5).PasteSpecial(-4122)
$ws.Cells.Item(106, 1).Value = 104
$ws.Cells.Item(106, 2).Value = 7749873
$ws.Cells.Item(106, 3).Value = "India Super League"
$ws.Cells.Item(106, 4).Value = "India Super League"
$ws.Cells.Item(106, 5).Value = 45359.45833333334
$ws.Cells.Item(106, 6).Value = "Jamshedpur FC"
$ws.Cells.Item(106, 7).Value = "Mumbai City FC"
$ws.Cells.Item(106, 11).Value = 3.75
$ws.Cells.Item(106, 12).Value = 3.5
$ws.Cells.Item(106, 13).Value = 1.95
$ws.Cells.Item(106, 14).Value = 3.75
$ws.Cells.Item(106, 15).Value = 3.75
$ws.Cells.Item(106, 16).Value = 1.909
$ws.Cells.Item(106, 17).Value = 0.5
$ws.Cells.Item(106, 18).Value = 1.9
$ws.Cells.Item(106, 19).Value = 1.9
$ws.Cells.Item(106, 20).Value = 2.75
$ws.Cells.Item(106, 21).Value = 1.85
$ws.Cells.Item(106, 22).Value = 1.95
$ws.Cells.Item(106, 23).Value = 0
$ws.Cells.Item(106, 24).Value = 0
$ws.Cells.Item(106, 25).Value = 0
$ws.Cells.Item(106, 26).Value = 0
$ws.Cells.Item(106, 27).Value = 0

# --- Row 107: new fixture ---
$ws.Range("A2").Copy()
$ws.Cells.Item(107, 1).PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Cells.Item(107, 5).PasteSpecial(-4122)
$ws.Cells.Item(107, 1).Value = 105
$ws.Cells.Item(107, 2).Value = 7751757
$ws.Cells.Item(107, 3).Value = "India Super League"
$ws.Cells.Item(107, 4).Value = "India Super League"
$ws.Cells.Item(107, 5).Value = 45360.45833333334
$ws.Cells.Item(107, 6).Value = "Chennaiyin FC"
$ws.Cells.Item(107, 7).Value = "Hyderabad FC"
$ws.Cells.Item(107, 11).Value = 1.333
$ws.Cells.Item(107, 12).Value = 4.5
$ws.Cells.Item(107, 13).Value = 8
$ws.Cells.Item(107, 14).Value = 1.333
$ws.Cells.Item(107, 15).Value = 4.5
$ws.Cells.Item(107, 16).Value = 8
$ws.Cells.Item(107, 17).Value = -1.5
$ws.Cells.Item(107, 18).Value = 2.025
$ws.Cells.Item(107, 19).Value = 1.775
$ws.Cells.Item(107, 20).Value = 2.75
$ws.Cells.Item(107, 21).Value = 1.9
$ws.Cells.Item(107, 22).Value = 1.9
$ws.Cells.Item(107, 23).Value = 0
$ws.Cells.Item(107, 24).Value = 0
$ws.Cells.Item(107, 25).Value = 0
$ws.Cells.Item(107, 26).Value = 0
$ws.Cells.Item(107, 27).Value = 0

$excel.CutCopyMode = $false

